$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.169.79"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "2.376.55"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.47"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.61"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.25%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -2.01%  "
$ws.Range("D9").Value = "2.377.74"
$ws.Range("E10").Value = "  +1.69%  "
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.348"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.06"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.80%  "
$ws.Range("D15").Value = "2.790.68"
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "61.070.86"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "2.378.16"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.92"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.14"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.20"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.70"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.34"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.71"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -10.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.49"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.16"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "509.23"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.13%  "
$ws.Range("D29").Value = "0.0₃0892"
$ws.Range("E29").Value = "  -4.21%  "
$ws.Range("E30").Value = "  +2.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.37"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.99%  "
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("E33").Value = "  -3.59%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.67"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.89"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.43"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.378"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.47"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "147.16"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.09%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "41.24"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "151.46"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +7.46%  "
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.59"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0521"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.26"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.576"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0907"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("E51").Value = "  +0.30%  "
